$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Grader increased points awarded for "addProduct() method" row (row 20)
# from 6 to 9, and updated the grading comment accordingly.
$ws.Range("E20").Value = 9
$ws.Range("F20").Value = "(-1) For not adding products if the customer exisis"

# Reflect the grader's last selection location.
$ws.Range("G24").Select()

$wb.Save()
